$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MLCO")

# Insert two new columns before column D (dates 2018-12-31 and 2018-09-30 quarter
# columns), shifting the existing quarterly data from D:K to F:M.
$ws.Range("D1:E1").EntireColumn.Insert()

# The inserted columns pick up the format of the column to their left (C) by
# default; re-apply the correct number formats (date format for the header
# row, number format for the data rows) by pulling them from the
# now-shifted former D/E columns (F/G) so the existing style indexes are
# reused instead of new ones being minted.
$ws.Range("F7:F102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$ws.Range("G7:G102").Copy()
$ws.Range("E7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the two new quarter columns with their values.
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 1396500
$ws.Range("E8").Value = 1220300
$ws.Range("D9").Value = 884800
$ws.Range("E9").Value = 814800
$ws.Range("D10").Value = 511700
$ws.Range("E10").Value = 405500
$ws.Range("D12").Value = 11300
$ws.Range("E12").Value = 4800
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 3200
$ws.Range("E14").Value = 200
$ws.Range("D15").Value = 149700
$ws.Range("E15").Value = 152900
$ws.Range("D17").Value = 1195700
$ws.Range("E17").Value = 1136800
$ws.Range("D18").Value = 200800
$ws.Range("E18").Value = 83500
$ws.Range("D20").Value = -2800
$ws.Range("E20").Value = -4100
$ws.Range("D21").Value = "NA"
$ws.Range("E21").Value = "NA"
$ws.Range("D22").Value = 74000
$ws.Range("E22").Value = 70800
$ws.Range("D23").Value = 124000
$ws.Range("E23").Value = 8500
$ws.Range("D24").Value = -6200
$ws.Range("E24").Value = 1300
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 130200
$ws.Range("E26").Value = 7200
$ws.Range("D27").Value = 128000
$ws.Range("E27").Value = 9600
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 2800
$ws.Range("E32").Value = 4100
$ws.Range("D33").Value = 128000
$ws.Range("E33").Value = 9600
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 128000
$ws.Range("E35").Value = 9600
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 1436600
$ws.Range("E41").Value = 1226500
$ws.Range("D42").Value = 91600
$ws.Range("E42").Value = 112800
$ws.Range("D43").Value = 249700
$ws.Range("E43").Value = 206000
$ws.Range("D44").Value = 40800
$ws.Range("E44").Value = 40100
$ws.Range("D45").Value = 138800
$ws.Range("E45").Value = 169700
$ws.Range("D46").Value = 1957500
$ws.Range("E46").Value = 1755200
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 5661700
$ws.Range("E48").Value = 5740400
$ws.Range("D49").Value = 1068600
$ws.Range("E49").Value = 1069700
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 189600
$ws.Range("E52").Value = 200900
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 8877400
$ws.Range("E54").Value = 8766200
$ws.Range("D57").Value = 24900
$ws.Range("E57").Value = 24300
$ws.Range("D58").Value = 430200
$ws.Range("E58").Value = 119500
$ws.Range("D59").Value = 1674900
$ws.Range("E59").Value = 1659000
$ws.Range("D60").Value = 2130000
$ws.Range("E60").Value = 1802800
$ws.Range("D61").Value = 3918700
$ws.Range("E61").Value = 3966800
$ws.Range("D62").Value = 82900
$ws.Range("E62").Value = 83200
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 6750000
$ws.Range("E66").Value = 6294800
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -703600
$ws.Range("E72").Value = -763800
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 2127300
$ws.Range("E76").Value = 2471400
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 128000
$ws.Range("E81").Value = 9600
$ws.Range("D83").Value = 0
$ws.Range("E83").Value = 0
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 0
$ws.Range("E89").Value = 0
$ws.Range("D91").Value = 0
$ws.Range("E91").Value = 0
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = 0
$ws.Range("E94").Value = 0
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 0
$ws.Range("E100").Value = 0
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 0
$ws.Range("E102").Value = 0
